$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the page-title heading near the top of the document.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# 2. At the end of the document, the final paragraph (the italic "Prompt: ..."
#    image-generation note) is replaced by two paragraphs:
#      - a new bold paragraph containing the page title text
#      - the same italic paragraph, but with its text swapped for the meta
#        description text that used to live at the top of the document
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Creature from the Black Lagoon for Free - Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the immersive gameplay and interactive features of Creature from the Black Lagoon slot game. Play for free and win big.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
